# Add "counting" sheet (Ark2) with parking-spot totals to the calculation
# workbook, and update the active-sheet/selection bookkeeping to match.

$wb = $excel.ActiveWorkbook

# --- Ark1: selection moves, it is no longer the active tab ---------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G17").Select()

# --- Ark2: new worksheet right after Ark1 ---------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Ark2"

$data = @(
    @("Sted", "Antall Plasser"),
    @("Inspiria", 125),
    @("Inspiria Bak", 40),
    @("Superland", 200),
    @("Quality Hotell", 205),
    @("Kiwi", 210),
    @("Politihuset", 85),
    @("Caverion", 45),
    @("K5 Bygget", 40),
    @("Tune Senter", 115),
    @("Adecco og If", 110),
    @("Fagforbundet", 110),
    @("Utenfor rundt Politiet", 85)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $data[$i][0]
    $ws2.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws2.Range("A14").Value = "Totalt"
$ws2.Range("B14").Formula = "=SUM(B2:B13)"

$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(2).AutoFit() | Out-Null

$ws2.Range("C16").Select()
